$d = $word.ActiveDocument

$pairs = @(
    @("2026-01-22 Thursday", "2026-01-23 Friday"),
    @("6+21=27", "78-0=78"),
    @("58-33=25", "58-55=3"),
    @("36-17=19", "91-24=67"),
    @("68-44=24", "49+10=59"),
    @("36-14=22", "17+28=45"),
    @("97-47=50", "22-12=10"),
    @("16+73=89", "29-19=10"),
    @("72+22=94", "81-10=71"),
    @("56+15=71", "26+29=55"),
    @("31+46=77", "13+4=17"),
    @("67-49=18", "6+77=83"),
    @("84-69=15", "49-25=24"),
    @("74+8=82", "56+5=61"),
    @("12+47=59", "75+8=83"),
    @("30-5=25", "77-29=48"),
    @("39+13=52", "62+35=97"),
    @("84-2=82", "56+6=62"),
    @("16+27=43", "10+78=88"),
    @("53+40=93", "40-4=36"),
    @("85-8=77", "53-47=6"),
    @("54+36=90", "18+22=40"),
    @("7+82=89", "20-11=9"),
    @("42+38=80", "70+20=90"),
    @("59-57=2", "44+42=86"),
    @("82-62=20", "36-6=30"),
    @("64-3=61", "12-6=6"),
    @("24+5=29", "9-3=6"),
    @("32+19=51", "20+37=57"),
    @("59-1=58", "30+45=75"),
    @("92-51=41", "62-29=33"),
    @("73-55=18", "65+26=91"),
    @("32+30=62", "33+32=65"),
    @("33+17=50", "98-47=51"),
    @("54+13=67", "7+33=40"),
    @("51-9=42", "8+57=65"),
    @("8+38=46", "58+0=58"),
    @("33+0=33", "76-46=30"),
    @("60+23=83", "83+0=83"),
    @("43+54=97", "38-26=12"),
    @("95-29=66", "17+64=81"),
    @("52+13=65", "51-17=34"),
    @("37-17=20", "3+22=25"),
    @("75-39=36", "79+9=88"),
    @("69+26=95", "14+67=81"),
    @("93-82=11", "46-30=16"),
    @("68-14=54", "17+45=62"),
    @("57-10=47", "87-8=79"),
    @("30+9=39", "8+4=12"),
    @("42-2=40", "16+8=24"),
    @("5+42=47", "4+45=49"),
    @("67-34=33", "68-11=57"),
    @("36+13=49", "6+81=87"),
    @("91+1=92", "87+5=92"),
    @("73-6=67", "67-24=43"),
    @("2+58=60", "25-24=1"),
    @("61-52=9", "19+76=95"),
    @("30+10=40", "10+32=42"),
    @("95+4=99", "45-29=16"),
    @("49-44=5", "23+42=65"),
    @("41-29=12", "95-10=85"),
    @("65+31=96", "59-2=57"),
    @("4+40=44", "28+48=76"),
    @("63+28=91", "34-4=30"),
    @("49+44=93", "37+3=40"),
    @("3+39=42", "65-39=26"),
    @("22+27=49", "97-61=36"),
    @("19-4=15", "26-23=3"),
    @("0+24=24", "25-18=7"),
    @("60-32=28", "39+21=60"),
    @("3-2=1", "0+90=90"),
    @("39-11=28", "46-28=18"),
    @("61+24=85", "15+71=86"),
    @("78-64=14", "26-12=14"),
    @("3+70=73", "45+36=81"),
    @("71-70=1", "87-82=5"),
    @("85-14=71", "66-50=16"),
    @("84-27=57", "44+39=83"),
    @("70-43=27", "56+11=67"),
    @("54+37=91", "82+0=82"),
    @("66+8=74", "19+67=86"),
    @("96+1=97", "87-33=54"),
    @("51+20=71", "62-45=17"),
    @("19+33=52", "22+75=97"),
    @("95-12=83", "36-12=24"),
    @("29+64=93", "14+6=20"),
    @("28-28=0", "19-8=11"),
    @("86-3=83", "57+10=67"),
    @("98-70=28", "24-17=7"),
    @("84+14=98", "18+69=87"),
    @("39+60=99", "55+22=77"),
    @("39+34=73", "11+76=87"),
    @("83-15=68", "83-74=9"),
    @("56-46=10", "52+11=63"),
    @("74-54=20", "44+37=81"),
    @("61-33=28", "26+33=59"),
    @("53-50=3", "57-30=27"),
    @("15+1=16", "73-3=70"),
    @("96-28=68", "79-8=71"),
    @("21+54=75", "1+15=16"),
    @("76+1=77", "30+5=35"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replacements complete: $($pairs.Count)"